# Applies the "Docs - minor update" commit to Functionality.docx
$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# ---------------------------------------------------------------------------
# 1) "brief guide below" -> "brief guide"
# ---------------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute(
    "we have prepared a brief guide below. Feel free to browse",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "we have prepared a brief guide. Feel free to browse", $wdReplaceOne)

# Word tracks the last edit position with a hidden "_GoBack" bookmark; it
# moves there the next time the document is opened. Re-anchor it at the end
# of the sentence we just touched (where the user's cursor last was).
$d.Bookmarks("_GoBack").Delete()

$markerRange = $d.Content
[void]$markerRange.Find.Execute(
    "have a great day!", $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, "", 0)
$markerRange.Collapse(0)
$markerRange.InsertAfter("@@GOBACKMARK@@")

$markerRange = $d.Content
[void]$markerRange.Find.Execute(
    "@@GOBACKMARK@@", $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange.Text = ""

Write-Output "done"
